# Apply the new table style to the table on slide 6 (the only table shape
# in the deck). PowerPoint table styles are identified by a GUID that maps
# to an entry either in ppt/tableStyles.xml or to one of PowerPoint's
# built-in table-style gallery entries; changing it is done with
# Table.ApplyStyle, not by assigning Table.Style directly.

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(6)
$shp = $s.Shapes.Item(2)
$tbl = $shp.Table

$tbl.ApplyStyle("{7F727C2F-BC66-473C-BC5F-4A0777269927}")
